# Updates cryptos list figures (price & 1h volume change) per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.337.37"
$ws.Range("E2").Value = "  +3.99%  "

$ws.Range("D3").Value = "1.594.08"
$ws.Range("E3").Value = "  +1.89%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "`'214.22"
$ws.Range("E5").Value = "  +1.67%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").Value = "`'23.98"
$ws.Range("E8").Value = "  +8.36%  "

$ws.Range("E9").Value = "  +0.84%  "

$ws.Range("E10").Value = "  +0.89%  "

$ws.Range("D11").Value = "`'0.0890"
$ws.Range("E11").Value = "  +2.13%  "

$ws.Range("D12").Value = "1.823.18"
$ws.Range("E12").Value = "  +2.00%  "

$ws.Range("D13").Value = "1.589.37"
$ws.Range("E13").Value = "  +1.52%  "

$ws.Range("D14").Value = "`'0.531"
$ws.Range("E14").Value = "  +2.60%  "

$ws.Range("E15").Value = "  -0.13%  "

$ws.Range("D16").Value = "28.359.61"
$ws.Range("E16").Value = "  +4.24%  "

$ws.Range("D17").Value = "`'63.12"

$ws.Range("D18").Value = "`'227.49"
$ws.Range("E18").Value = "  +4.38%  "

$ws.Range("D19").Value = "0.0₃0710"
$ws.Range("E19").Value = "  +1.16%  "

$ws.Range("D20").Value = "`'7.48"
$ws.Range("E20").Value = "  +0.34%  "

$ws.Range("E21").Value = "  -0.02%  "

$ws.Range("E22").Value = "  -0.70%  "

$ws.Range("D23").Value = "`'9.33"
$ws.Range("E23").Value = "  -0.51%  "

$ws.Range("E24").Value = "  +0.70%  "

$ws.Range("D25").Value = "`'151.59"
$ws.Range("E25").Value = "  -0.08%  "

$ws.Range("D26").Value = "`'15.20"

$ws.Range("E28").Value = "  -0.51%  "

$ws.Range("E29").Value = "  +0.09%  "

$ws.Range("E30").Value = "  +0.63%  "

$ws.Range("D31").Value = "`'0.0475"
$ws.Range("E31").Value = "  +1.25%  "

$ws.Range("D32").Value = "`'3.23"
$ws.Range("E32").Value = "  -0.06%  "

$ws.Range("E33").Value = "  -0.98%  "

$ws.Range("D34").Value = "1.398.09"
$ws.Range("E34").Value = "  -3.97%  "

$ws.Range("E35").Value = "  -1.54%  "

$ws.Range("E36").Value = "  -5.94%  "

$ws.Range("E37").Value = "  +0.21%  "

$ws.Range("E38").Value = "  +0.63%  "

$ws.Range("E39").Value = "  +8.74%  "

$ws.Range("E40").Value = "  +0.35%  "

$ws.Range("E41").Value = "  -0.08%  "

$ws.Range("D42").Value = "`'5.72"
$ws.Range("E42").Value = "  -2.51%  "

$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").Value = "`'1.00"
$ws.Range("E43").Value = "  -0.02%  "

$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "`'1.89"
$ws.Range("E44").Value = "  +7.72%  "

$ws.Range("D45").Value = "`'0.984"
$ws.Range("E45").Value = "  +0.21%  "

$ws.Range("D46").Value = "`'64.41"
$ws.Range("E46").Value = "  -0.05%  "

$ws.Range("D47").Value = "1.733.66"
$ws.Range("E47").Value = "  +2.03%  "

$ws.Range("B48").Value = "mCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$ws.Range("D48").Value = "`'2.14"
$ws.Range("E48").Value = "  +0.40%  "

$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "`'87.52"
$ws.Range("E49").Value = "  +1.96%  "

$ws.Range("E50").Value = "  +0.04%  "

$ws.Range("D51").Value = "`'0.0527"
$ws.Range("E51").Value = "  +0.33%  "
